$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.1.0 -> 1.1.1
$ws.Range("B3").Value = "1.1.1"

# Experimental: (empty) -> "false"
# Assigning the literal string "false" directly would be auto-coerced to the
# Boolean value FALSE by Excel's type inference, so instead write it as a
# formula producing the text "false" and then convert it to a literal value
# via copy / paste-special (values only), which preserves the Text type.
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Date: 2022-05-04T10:16:52-05:00 -> 2022-10-21T09:04:31-05:00
$ws.Range("B8").Value = "2022-10-21T09:04:31-05:00"
